# Bump the "Förändrad" date (column C) by one day for all data rows (2-297)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 297
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
